$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.396007019381747
$ws.Range("C2").Value = 2.396007019381743
$ws.Range("D2").Value = 2.396007019381743

$ws.Range("B3").Value = 0.03343662252085296
$ws.Range("C3").Value = 0.03854421396097765
$ws.Range("D3").Value = 0.1751469250153707

$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03729126592179306
$ws.Range("C4").Value = 0.03218526624368307
$ws.Range("D4").Value = 0.04626265997781376

$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.03295216155382702
$ws.Range("C5").Value = 0.03953824054679839
$ws.Range("D5").Value = 0.06120857651209916
